# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps for the 55851552-7366-4db0-a33e-6e891440e385 / 7bba6331-cfdc-4895-adab-6ebd8a90c20e
# entries (row 3, which shares its text with row 4) on both the zh-cn and
# de-de sheets, reflecting a newer report generation run.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E3").Value = "2016-03-11 10:24:49"
$ws_zhcn.Range("H3").Value = "2016-03-11 10:25:07"
$ws_zhcn.Range("E4").Value = "2016-03-11 10:24:49"
$ws_zhcn.Range("H4").Value = "2016-03-11 10:25:07"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E3").Value = "2016-03-11 10:24:52"
$ws_dede.Range("H3").Value = "2016-03-11 10:25:12"
$ws_dede.Range("E4").Value = "2016-03-11 10:24:52"
$ws_dede.Range("H4").Value = "2016-03-11 10:25:12"
